$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells so numeric-looking
# strings like "1.00" or "208.34" are preserved as text, matching
# the original inline-string cell type instead of being parsed as numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.021.90'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.563.09'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('E4').Value = '  +0.43%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.34'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('E6').Value = '  +0.58%  '
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0599'
$ws.Range('E10').Value = '  +1.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0856'
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.783.48'
$ws.Range('E12').Value = '  +0.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.533.33'
$ws.Range('E13').Value = '  -1.44%  '
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.522'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.011.98'
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.96'
$ws.Range('E17').Value = '  +0.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0708'
$ws.Range('E18').Value = '  +1.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '216.48'
$ws.Range('E19').Value = '  -0.78%  '
$ws.Range('E20').Value = '  +0.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.01'
$ws.Range('E21').Value = '  +0.36%  '
$ws.Range('E22').Value = '  +1.73%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.01'
$ws.Range('E25').Value = '  -0.86%  '
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.13'
$ws.Range('E27').Value = '  +1.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.106'
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.01'
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('E30').Value = '  +1.19%  '
$ws.Range('E31').Value = '  +2.68%  '
$ws.Range('E32').Value = '  -0.23%  '
$ws.Range('E33').Value = '  +2.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.427.60'
$ws.Range('E34').Value = '  +0.29%  '
$ws.Range('E35').Value = '  +1.37%  '
$ws.Range('E36').Value = '  +8.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.34'
$ws.Range('E37').Value = '  +2.45%  '
$ws.Range('E38').Value = '  +0.42%  '
$ws.Range('E39').Value = '  +2.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.86'
$ws.Range('E40').Value = '  +2.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.809'
$ws.Range('E41').Value = '  -0.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.01'
$ws.Range('E42').Value = '  +0.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +1.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.32'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.83'
$ws.Range('E45').Value = '  +0.49%  '
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.698.48'
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.20'
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('E49').Value = '  +5.71%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0519'
$ws.Range('E50').Value = '  -0.21%  '
$ws.Range('E51').Value = '  +0.36%  '

Write-Host "Applied 76 cell updates"
